$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 8.307
$ws.Cells.Item(3, 1).Value = -21.8085
$ws.Cells.Item(3, 4).Value = -7.414599999999996
$ws.Cells.Item(6, 5).Value = 16.54880000000001
$ws.Cells.Item(12, 4).Value = -7.265700000000002
$ws.Cells.Item(14, 1).Value = -21.86049999999999
$ws.Cells.Item(19, 5).Value = 16.26209999999999
$ws.Cells.Item(21, 1).Value = -20.15989999999999
$ws.Cells.Item(23, 1).Value = -20.12719999999998
$ws.Cells.Item(24, 4).Value = -7.370800000000004
$ws.Cells.Item(24, 5).Value = 16.96730000000001
$ws.Cells.Item(25, 1).Value = -21.8689
$ws.Cells.Item(25, 2).Value = 5.850400000000002
$ws.Cells.Item(25, 4).Value = -8.94579999999999
$ws.Cells.Item(26, 1).Value = -21.08039999999996
$ws.Cells.Item(27, 2).Value = 6.322300000000001
$ws.Cells.Item(29, 1).Value = -20.62419999999998
$ws.Cells.Item(30, 5).Value = 15.654
$ws.Cells.Item(31, 2).Value = 4.903199999999999
$ws.Cells.Item(31, 5).Value = 16.232
$ws.Cells.Item(33, 5).Value = 16.86280000000002
$ws.Cells.Item(39, 2).Value = 9.566300000000002
$ws.Cells.Item(42, 5).Value = 16.39259999999999
$ws.Cells.Item(48, 2).Value = 5.096000000000001
$ws.Cells.Item(50, 4).Value = -8.195500000000001
$ws.Cells.Item(51, 2).Value = 5.573099999999999
$ws.Cells.Item(52, 2).Value = 5.143599999999999
$ws.Cells.Item(53, 1).Value = -22.21700000000001
$ws.Cells.Item(53, 4).Value = -6.0812
$ws.Cells.Item(55, 2).Value = 5.924999999999998
$ws.Cells.Item(55, 5).Value = 16.4927
$ws.Cells.Item(56, 2).Value = 5.368499999999996
$ws.Cells.Item(57, 1).Value = -22.1401
$ws.Cells.Item(57, 2).Value = 4.592999999999996
$ws.Cells.Item(57, 4).Value = -8.478300000000004
$ws.Cells.Item(58, 5).Value = 16.73840000000001
$ws.Cells.Item(59, 1).Value = -22.2893
$ws.Cells.Item(61, 4).Value = -7.678299999999999
$ws.Cells.Item(63, 4).Value = -7.756400000000005
$ws.Cells.Item(65, 5).Value = 16.72160000000001
$ws.Cells.Item(69, 1).Value = -21.58579999999999
$ws.Cells.Item(70, 4).Value = -8.1317
$ws.Cells.Item(70, 5).Value = 16.90799999999999
$ws.Cells.Item(73, 2).Value = 8.602899999999998
$ws.Cells.Item(75, 5).Value = 16.58440000000001
$ws.Cells.Item(79, 1).Value = -20.3634
$ws.Cells.Item(83, 1).Value = -22.04249999999999
$ws.Cells.Item(83, 5).Value = 16.66520000000001
$ws.Cells.Item(86, 4).Value = -8.741400000000006
$ws.Cells.Item(86, 5).Value = 15.9304
$ws.Cells.Item(89, 2).Value = 4.799799999999994
$ws.Cells.Item(90, 2).Value = 5.811700000000005
$ws.Cells.Item(91, 1).Value = -21.37130000000002
$ws.Cells.Item(92, 2).Value = 4.830799999999996
$ws.Cells.Item(93, 1).Value = -20.79469999999998
$ws.Cells.Item(96, 5).Value = 15.97359999999999
$ws.Cells.Item(97, 5).Value = 17.04170000000002
$ws.Cells.Item(98, 4).Value = -8.739499999999998
$ws.Cells.Item(100, 4).Value = -8.643899999999999
$ws.Cells.Item(102, 4).Value = -7.815399999999994
